# Update the workbook window chrome (xWindow/windowWidth in bookViews).
$aw = $excel.ActiveWindow
$aw.Left = 9180
$aw.Width = 19620
$aw.Top = 1180
$aw.Height = 15520

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (http://www.1ne.org.uk/) is removed entirely - clear it so the row
# disappears from sheetData and the used range starts at row 2.
$ws.Range("A1").ClearContents()

# Rows 2-13 get replaced with a new list of charity URLs (written in the
# same order the source data was authored in, so the shared-string table
# comes out in the same sequence as the target workbook).
$ws.Range("A8").Value = "http://www.youthtalk.org.uk/"
$ws.Range("A7").Value = "http://www.womensaidnel.org/"
$ws.Range("A6").Value = "http://www.swingsandsmiles.co.uk/"
$ws.Range("A5").Value = "http://www.southendcarers.co.uk/"
$ws.Range("A4").Value = "http://www.psspeople.com/"
$ws.Range("A3").Value = "http://www.no-secrets.org.uk/"
$ws.Range("A2").Value = "http://www.mindincroydon.org.uk/"
$ws.Range("A13").Value = "http://www.middlesbroughandstocktonmind.org.uk/"
$ws.Range("A12").Value = "https://www.ymcadlg.org/"
$ws.Range("A11").Value = "https://www.place2be.org.uk/"
$ws.Range("A10").Value = "https://www.kidscape.org.uk/"
$ws.Range("A9").Value = "https://www.disc-vol.org.uk/"

# Rows 14-15 (old tail of the list) are removed entirely.
$ws.Range("A14:A15").ClearContents()

# Update the selection to match the new used range tail.
$ws.Range("A9:A13").Select()
